# Weekly data refresh: a new observation is published for
# "Feria Lagunitas de Puerto Montt - Apio" at the top of the data block
# (row 279), pushing all existing rows (279-424) down by one
# (to 280-425).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 279; this shifts rows 279:424 -> 280:425
# and grows the sheet from A1:R424 to A1:R425.
$ws.Rows("279:279").Insert()

# Populate the new row with the latest observation.
$ws.Range("A279").Value = 4
$ws.Range("B279").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C279").Value = "Los Lagos"
$ws.Range("D279").Value = 45001
$ws.Range("E279").Value = 10
$ws.Range("F279").Value = 100112017
$ws.Range("G279").Value = "Apio"
$ws.Range("H279").Value = "Americana (o)"
$ws.Range("I279").Value = "Primera"
$ws.Range("J279").Value = 25
$ws.Range("K279").Value = 12000
$ws.Range("L279").Value = 12000
$ws.Range("M279").Value = 12000
$ws.Range("N279").Value = "$/docena de matas"
$ws.Range("O279").Value = "Región de Coquimbo"
$ws.Range("P279").Value = 2000
$ws.Range("Q279").Value = 6
$ws.Range("R279").Value = "Hortaliza"
